# Weekly refresh of Hortaliza / Poroto granado price rows (Feria Lagunitas de Puerto Montt).
# Each data row (2-25) is updated in place: Fecha, Volumen, Precio minimo/maximo/promedio,
# Origen and Precio $/Kg move to the values from the newly-reported week; Calidad (I) only
# changes where the reported quality tier itself changed (rows 16 and 20).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44204
$ws.Range("J2").Value = 50
$ws.Range("K2").Value = 37000
$ws.Range("L2").Value = 37000
$ws.Range("M2").Value = 37000
$ws.Range("P2").Value = 1480

# Row 3
$ws.Range("D3").Value = 44575
$ws.Range("J3").Value = 80
$ws.Range("K3").Value = 35000
$ws.Range("L3").Value = 35000
$ws.Range("M3").Value = 35000
$ws.Range("O3").Value = "Región Metropolitana"
$ws.Range("P3").Value = 1400

# Row 4
$ws.Range("D4").Value = 44243
$ws.Range("J4").Value = 50
$ws.Range("K4").Value = 33000
$ws.Range("L4").Value = 33000
$ws.Range("M4").Value = 33000
$ws.Range("P4").Value = 1320

# Row 6
$ws.Range("D6").Value = 44201
$ws.Range("J6").Value = 60
$ws.Range("K6").Value = 30000
$ws.Range("L6").Value = 30000
$ws.Range("M6").Value = 30000
$ws.Range("P6").Value = 1200

# Row 7
$ws.Range("D7").Value = 44271
$ws.Range("J7").Value = 40
$ws.Range("K7").Value = 30000
$ws.Range("L7").Value = 30000
$ws.Range("M7").Value = 30000
$ws.Range("O7").Value = "Región del Maule"
$ws.Range("P7").Value = 1200

# Row 8
$ws.Range("D8").Value = 44239
$ws.Range("J8").Value = 60
$ws.Range("K8").Value = 35000
$ws.Range("L8").Value = 35000
$ws.Range("M8").Value = 35000
$ws.Range("P8").Value = 1400

# Row 9
$ws.Range("D9").Value = 44236
$ws.Range("J9").Value = 40
$ws.Range("K9").Value = 38000
$ws.Range("L9").Value = 38000
$ws.Range("M9").Value = 38000
$ws.Range("P9").Value = 1520

# Row 10
$ws.Range("D10").Value = 44222
$ws.Range("K10").Value = 40000
$ws.Range("L10").Value = 40000
$ws.Range("M10").Value = 40000
$ws.Range("P10").Value = 1600

# Row 11
$ws.Range("D11").Value = 44215
$ws.Range("J11").Value = 60
$ws.Range("O11").Value = "Región del Maule"

# Row 14
$ws.Range("D14").Value = 44232
$ws.Range("J14").Value = 40
$ws.Range("K14").Value = 40000
$ws.Range("L14").Value = 40000
$ws.Range("M14").Value = 40000
$ws.Range("P14").Value = 1600

# Row 15
$ws.Range("D15").Value = 44582
$ws.Range("K15").Value = 35000
$ws.Range("L15").Value = 35000
$ws.Range("M15").Value = 35000
$ws.Range("O15").Value = "Región Metropolitana"
$ws.Range("P15").Value = 1400

# Row 16
$ws.Range("D16").Value = 44582
$ws.Range("I16").Value = "Segunda"
$ws.Range("J16").Value = 40
$ws.Range("K16").Value = 27000
$ws.Range("L16").Value = 27000
$ws.Range("M16").Value = 27000
$ws.Range("O16").Value = "Región Metropolitana"
$ws.Range("P16").Value = 1080

# Row 17
$ws.Range("D17").Value = 44572
$ws.Range("J17").Value = 80
$ws.Range("O17").Value = "Región Metropolitana"

# Row 18
$ws.Range("D18").Value = 44246
$ws.Range("K18").Value = 31000
$ws.Range("L18").Value = 31000
$ws.Range("M18").Value = 31000
$ws.Range("P18").Value = 1240

# Row 19
$ws.Range("D19").Value = 44218
$ws.Range("J19").Value = 60
$ws.Range("K19").Value = 42000
$ws.Range("L19").Value = 42000
$ws.Range("M19").Value = 42000
$ws.Range("O19").Value = "Región del Maule"
$ws.Range("P19").Value = 1680

# Row 20
$ws.Range("D20").Value = 44211
$ws.Range("I20").Value = "Primera"
$ws.Range("K20").Value = 42000
$ws.Range("L20").Value = 42000
$ws.Range("M20").Value = 42000
$ws.Range("O20").Value = "Región del Maule"
$ws.Range("P20").Value = 1680

# Row 21
$ws.Range("D21").Value = 44203
$ws.Range("J21").Value = 20
$ws.Range("K21").Value = 30000
$ws.Range("L21").Value = 30000
$ws.Range("M21").Value = 30000
$ws.Range("P21").Value = 1200

# Row 22
$ws.Range("D22").Value = 44253
$ws.Range("J22").Value = 80
$ws.Range("K22").Value = 30000
$ws.Range("L22").Value = 30000
$ws.Range("M22").Value = 30000
$ws.Range("P22").Value = 1200

# Row 23
$ws.Range("D23").Value = 44202
$ws.Range("J23").Value = 30
$ws.Range("K23").Value = 30000
$ws.Range("L23").Value = 30000
$ws.Range("M23").Value = 30000
$ws.Range("P23").Value = 1200

# Row 24
$ws.Range("D24").Value = 44225
$ws.Range("J24").Value = 60
$ws.Range("K24").Value = 32000
$ws.Range("L24").Value = 32000
$ws.Range("M24").Value = 32000
$ws.Range("P24").Value = 1280

# Row 25
$ws.Range("D25").Value = 44250
$ws.Range("J25").Value = 70
$ws.Range("K25").Value = 30000
$ws.Range("L25").Value = 30000
$ws.Range("M25").Value = 30000
$ws.Range("P25").Value = 1200
